$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.337.46'
$ws.Range("E2").Value = '  -3.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.499.37'
$ws.Range("E3").Value = '  -4.68%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.93'
$ws.Range("E5").Value = '  -2.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.44'
$ws.Range("E6").Value = '  -6.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.498.41'
$ws.Range("E7").Value = '  -4.62%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("E9").Value = '  -3.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.139'
$ws.Range("E10").Value = '  -3.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.94'
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -4.21%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000219'
$ws.Range("E13").Value = '  -4.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.091.53'
$ws.Range("E14").Value = '  -4.62%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '31.45'
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.496.53'
$ws.Range("E16").Value = '  -4.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.253.99'
$ws.Range("E17").Value = '  -3.28%  '
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("E19").Value = '  -2.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.98'
$ws.Range("E20").Value = '  -5.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '445.97'
$ws.Range("E21").Value = '  -4.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.97'
$ws.Range("E22").Value = '  -12.73%  '
$ws.Range("E23").Value = '  -4.60%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.21'
$ws.Range("E24").Value = '  -2.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000130'
$ws.Range("E25").Value = '  +6.09%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.639.77'
$ws.Range("E27").Value = '  -4.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  -7.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.25'
$ws.Range("E29").Value = '  -5.16%  '
$ws.Range("E30").Value = '  -4.73%  '
$ws.Range("E31").Value = '  -6.88%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.67'
$ws.Range("E34").Value = '  -3.45%  '
$ws.Range("E35").Value = '  -4.08%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.488.30'
$ws.Range("E36").Value = '  -4.99%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.84'
$ws.Range("E37").Value = '  -6.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.96'
$ws.Range("E38").Value = '  -3.71%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("E41").Value = '  -0.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '173.86'
$ws.Range("E42").Value = '  -2.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0875'
$ws.Range("E43").Value = '  -1.75%  '
$ws.Range("E44").Value = '  -6.35%  '
$ws.Range("E45").Value = '  -4.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.41'
$ws.Range("E46").Value = '  -2.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.84'
$ws.Range("E47").Value = '  -4.56%  '
$ws.Range("E48").Value = '  +6.24%  '
$ws.Range("E49").Value = '  -5.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.54'
$ws.Range("E50").Value = '  -4.06%  '
$ws.Range("E51").Value = '  -3.69%  '
